$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D, E, G keep their original text (inline-string style) representation
# by forcing Text format before assigning the new values, so Excel does not
# auto-convert numeric-looking / percentage-looking strings into numbers.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "329.93"
$ws.Range("E2").Value = "7.46%"
$ws.Range("G2").Value = "9"

# Row 3
$ws.Range("D3").Value = "39.96"
$ws.Range("E3").Value = "7.47%"
$ws.Range("G3").Value = "9"

# Row 4
$ws.Range("D4").Value = "5.268"
$ws.Range("E4").Value = "2.85%"
$ws.Range("G4").Value = "9"

# Row 5
$ws.Range("D5").Value = "0.08100"
$ws.Range("E5").Value = "4.04%"
$ws.Range("G5").Value = "9"

# Row 6
$ws.Range("D6").Value = "4.522"
$ws.Range("E6").Value = "1.94%"
$ws.Range("G6").Value = "9"

# Row 7
$ws.Range("D7").Value = "8.611"
$ws.Range("E7").Value = "4.44%"
$ws.Range("G7").Value = "9"

# Row 8
$ws.Range("D8").Value = "1.928"
$ws.Range("E8").Value = "2.51%"
$ws.Range("G8").Value = "9"

# Row 9
$ws.Range("D9").Value = "2.977"
$ws.Range("E9").Value = "-0.50%"
$ws.Range("G9").Value = "9"

# Row 10
$ws.Range("D10").Value = "0.9350"
$ws.Range("E10").Value = "0.86%"
$ws.Range("G10").Value = "9"

# Row 11
$ws.Range("D11").Value = "0.1323"
$ws.Range("E11").Value = "21.58%"
$ws.Range("G11").Value = "9"

# Row 12
$ws.Range("D12").Value = "0.1972"
$ws.Range("E12").Value = "3.23%"
$ws.Range("G12").Value = "9"

# Row 13
$ws.Range("D13").Value = "0.09265"
$ws.Range("E13").Value = "3.77%"
$ws.Range("G13").Value = "9"

# Row 14
$ws.Range("D14").Value = "0.03565"
$ws.Range("E14").Value = "6.10%"
$ws.Range("G14").Value = "9"

# Row 15
$ws.Range("D15").Value = "0.09575"
$ws.Range("E15").Value = "-0.09%"
$ws.Range("G15").Value = "9"

# Row 16
$ws.Range("D16").Value = "0.001334"
$ws.Range("E16").Value = "-3.10%"
$ws.Range("G16").Value = "9"

# Row 17
$ws.Range("D17").Value = "0.006137"
$ws.Range("E17").Value = "6.65%"
$ws.Range("G17").Value = "9"

# Row 18
$ws.Range("D18").Value = "3.371"
$ws.Range("E18").Value = "-4.67%"
$ws.Range("G18").Value = "9"

# Row 19
$ws.Range("D19").Value = "0.3524"
$ws.Range("E19").Value = "4.62%"
$ws.Range("G19").Value = "9"

# Row 20
$ws.Range("D20").Value = "6.986"
$ws.Range("E20").Value = "10.77%"
$ws.Range("G20").Value = "9"

# Row 21
$ws.Range("D21").Value = "0.1322"
$ws.Range("E21").Value = "3.79%"
$ws.Range("G21").Value = "9"

# Row 22
$ws.Range("D22").Value = "0.2561"
$ws.Range("E22").Value = "2.14%"
$ws.Range("G22").Value = "9"

# Row 23
$ws.Range("D23").Value = "0.04425"
$ws.Range("E23").Value = "1.11%"
$ws.Range("G23").Value = "9"

# Row 24
$ws.Range("D24").Value = "0.001223"
$ws.Range("E24").Value = "2.49%"
$ws.Range("G24").Value = "9"

# Row 25
$ws.Range("D25").Value = "0.004312"
$ws.Range("E25").Value = "1.59%"
$ws.Range("G25").Value = "9"

# Row 26
$ws.Range("E26").Value = "-8.70%"
$ws.Range("G26").Value = "9"

# Row 27
$ws.Range("D27").Value = "0.0003992"
$ws.Range("E27").Value = "0.00%"
$ws.Range("G27").Value = "9"

# Row 28
$ws.Range("G28").Value = "9"

# Row 29
$ws.Range("G29").Value = "9"

# Row 30
$ws.Range("G30").Value = "9"

# Row 31
$ws.Range("G31").Value = "9"

# Row 32
$ws.Range("G32").Value = "9"

# Row 33
$ws.Range("G33").Value = "9"

# Row 34
$ws.Range("G34").Value = "9"

# Row 35
$ws.Range("G35").Value = "9"

# Row 36
$ws.Range("G36").Value = "9"

# Row 37
$ws.Range("G37").Value = "9"

# Row 38
$ws.Range("G38").Value = "9"

# Row 39
$ws.Range("D39").Value = "0.02512"
$ws.Range("E39").Value = "16.95%"
$ws.Range("G39").Value = "9"

# Row 40
$ws.Range("D40").Value = "0.05156"
$ws.Range("E40").Value = "2.64%"
$ws.Range("G40").Value = "9"

# Row 41
$ws.Range("D41").Value = "0.007670"
$ws.Range("E41").Value = "2.80%"
$ws.Range("G41").Value = "9"

# Row 42
$ws.Range("G42").Value = "9"

# Row 43
$ws.Range("D43").Value = "0.009193"
$ws.Range("E43").Value = "5.66%"
$ws.Range("G43").Value = "9"

# Row 44
$ws.Range("D44").Value = "0.002172"
$ws.Range("E44").Value = "2.61%"
$ws.Range("G44").Value = "9"

# Row 45
$ws.Range("D45").Value = "0.01025"
$ws.Range("E45").Value = "28.18%"
$ws.Range("G45").Value = "9"

# Row 46
$ws.Range("D46").Value = "0.00006654"
$ws.Range("E46").Value = "1.23%"
$ws.Range("G46").Value = "9"

# Row 47
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "-0.26%"
$ws.Range("G47").Value = "9"

# Row 48
$ws.Range("E48").Value = "147.45%"
$ws.Range("G48").Value = "9"

# Row 49
$ws.Range("E49").Value = "1.61%"
$ws.Range("G49").Value = "9"

# Row 50
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "-0.26%"
$ws.Range("G50").Value = "9"

# Row 51
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "-0.26%"
$ws.Range("G51").Value = "9"
